$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows at 109-112 for the new week (shifts old rows 109-133 down to 113-137)
$ws.Rows("109:112").Insert()

# Constant columns for every data row on this sheet
$colA = 1
$colB = 'Agrícola del Norte S.A. de Arica'
$colC = 'Arica y Parinacota'
$colE = 15
$colF = 100114001
$colG = 'Papa'
$colQ = 25
$colR = 'Hortaliza'

# Row 109
$ws.Cells.Item(109, 1).Value = $colA
$ws.Cells.Item(109, 2).Value = $colB
$ws.Cells.Item(109, 3).Value = $colC
$ws.Cells.Item(109, 4).Value = 44943
$ws.Cells.Item(109, 5).Value = $colE
$ws.Cells.Item(109, 6).Value = $colF
$ws.Cells.Item(109, 7).Value = $colG
$ws.Cells.Item(109, 8).Value = 'Asterix'
$ws.Cells.Item(109, 9).Value = '1a (cosecha)'
$ws.Cells.Item(109, 10).Value = 1000
$ws.Cells.Item(109, 11).Value = 15000
$ws.Cells.Item(109, 12).Value = 16000
$ws.Cells.Item(109, 13).Value = 15500
$ws.Cells.Item(109, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(109, 15).Value = 'Provincia de Melipilla'
$ws.Cells.Item(109, 16).Value = 620
$ws.Cells.Item(109, 17).Value = $colQ
$ws.Cells.Item(109, 18).Value = $colR

# Row 110
$ws.Cells.Item(110, 1).Value = $colA
$ws.Cells.Item(110, 2).Value = $colB
$ws.Cells.Item(110, 3).Value = $colC
$ws.Cells.Item(110, 4).Value = 44943
$ws.Cells.Item(110, 5).Value = $colE
$ws.Cells.Item(110, 6).Value = $colF
$ws.Cells.Item(110, 7).Value = $colG
$ws.Cells.Item(110, 8).Value = 'Asterix'
$ws.Cells.Item(110, 9).Value = '1a (cosecha)'
$ws.Cells.Item(110, 10).Value = 1300
$ws.Cells.Item(110, 11).Value = 14000
$ws.Cells.Item(110, 12).Value = 15000
$ws.Cells.Item(110, 13).Value = 14462
$ws.Cells.Item(110, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(110, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(110, 16).Value = 578
$ws.Cells.Item(110, 17).Value = $colQ
$ws.Cells.Item(110, 18).Value = $colR

# Row 111
$ws.Cells.Item(111, 1).Value = $colA
$ws.Cells.Item(111, 2).Value = $colB
$ws.Cells.Item(111, 3).Value = $colC
$ws.Cells.Item(111, 4).Value = 44943
$ws.Cells.Item(111, 5).Value = $colE
$ws.Cells.Item(111, 6).Value = $colF
$ws.Cells.Item(111, 7).Value = $colG
$ws.Cells.Item(111, 8).Value = 'Red Lady'
$ws.Cells.Item(111, 9).Value = '1a (cosecha)'
$ws.Cells.Item(111, 10).Value = 1100
$ws.Cells.Item(111, 11).Value = 14000
$ws.Cells.Item(111, 12).Value = 16000
$ws.Cells.Item(111, 13).Value = 15182
$ws.Cells.Item(111, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(111, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(111, 16).Value = 607
$ws.Cells.Item(111, 17).Value = $colQ
$ws.Cells.Item(111, 18).Value = $colR

# Row 112
$ws.Cells.Item(112, 1).Value = $colA
$ws.Cells.Item(112, 2).Value = $colB
$ws.Cells.Item(112, 3).Value = $colC
$ws.Cells.Item(112, 4).Value = 44943
$ws.Cells.Item(112, 5).Value = $colE
$ws.Cells.Item(112, 6).Value = $colF
$ws.Cells.Item(112, 7).Value = $colG
$ws.Cells.Item(112, 8).Value = 'Rosara'
$ws.Cells.Item(112, 9).Value = '1a (cosecha)'
$ws.Cells.Item(112, 10).Value = 1000
$ws.Cells.Item(112, 11).Value = 14000
$ws.Cells.Item(112, 12).Value = 15000
$ws.Cells.Item(112, 13).Value = 14500
$ws.Cells.Item(112, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(112, 15).Value = 'Región del Maule'
$ws.Cells.Item(112, 16).Value = 580
$ws.Cells.Item(112, 17).Value = $colQ
$ws.Cells.Item(112, 18).Value = $colR
